# The deck's Design (slide master) theme is switched from the custom
# "Integral" color scheme to the built-in "Office" color scheme - i.e.
# the 12 theme colors (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# that live in ppt/theme/theme1.xml are replaced with the stock Office
# theme values, as if a different theme/variant had been applied from
# the Design tab.

$p   = $ppt.ActivePresentation
$sm  = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# Values are the standard VBA BGR-packed long for each "Office" theme
# srgbClr hex value:
#   dk1=000000 lt1=FFFFFF dk2=44546A lt2=E7E6E6
#   accent1=5B9BD5 accent2=ED7D31 accent3=A5A5A5 accent4=FFC000
#   accent5=4472C4 accent6=70AD47 hlink=0563C1 folHlink=954F72

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
